# Results_CNN sheet: add "Scheduler" / "Min. LR" columns (U, V) and append
# one new training-run row (row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: U1 "Scheduler", V1 "Min. LR" ------------------------
$ws.Range("U1").Value = "Scheduler"
$ws.Range("V1").Value = "Min. LR"
# Match the existing header styling (bold font, border, centered/top align)
# by copying the format from the last header cell (T1) onto the new ones.
$ws.Range("T1").Copy()
$ws.Range("U1:V1").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 17 previously had "None" parked in the "Learning rate decay"
#     (S) column; with the new Scheduler column that value belongs in row
#     18 instead, so S17 is cleared back to blank. ---------------------------
$ws.Range("S17").Value = ""

# --- New data row 18 ---------------------------------------------------------
$ws.Range("A18").Value = "2024-1-5 11:1:58"
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 64
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = "ADAM"
$ws.Range("F18").Value = "CEL"
$ws.Range("G18").Value = 35.9
$ws.Range("H18").Value = 32
$ws.Range("I18").Value = 1.9151
$ws.Range("J18").Value = 1.7357
$ws.Range("K18").Value = 26.0162
$ws.Range("L18").Value = "FER2013"
$ws.Range("M18").Value = "cpu"
$ws.Range("N18").Value = 4
$ws.Range("O18").Value = 2
$ws.Range("P18").Value = "Stationær"
$ws.Range("Q18").Value = 71.7
$ws.Range("R18").Value = 0.005
$ws.Range("S18").Value = ""
$ws.Range("T18").Value = 0.1
$ws.Range("U18").Value = "None"
$ws.Range("V18").Value = 0
